$wb = $excel.ActiveWorkbook

$wsDatos = $wb.Worksheets.Item("Datos")
$wsErrores = $wb.Worksheets.Item("Errores")

# --- Move existing rows 13 & 15 down to 15 & 17 on "Errores" sheet, and
# --- insert the two new comment rows at 12 and 13.

# Capture current values before moving anything.
$b13 = $wsErrores.Range("B13").Value2
$c13 = $wsErrores.Range("C13").Value2
$b15 = $wsErrores.Range("B15").Value2
$c15 = $wsErrores.Range("C15").Value2

# Move row 15 -> row 17 first (so we don't overwrite anything).
$wsErrores.Range("B17").Value = $b15
$wsErrores.Range("C17").Value = $c15
$wsErrores.Range("B15").Value = $b13
$wsErrores.Range("C15").Value = $c13

# Clear old row 13 content, then write the two new rows (12 & 13).
$wsErrores.Range("B13").ClearContents()
$wsErrores.Range("C12").Value = "Muchos de los metodos no usan la implementacion del toBO por tanto se pueden generar muchos errores al momento de mapear la info"
$wsErrores.Range("C13").Value = "El metodo de cerrar subasta no esta persistiendo los cambios"

# --- Page setup on "Errores" sheet (adds pageSetup element / printer settings)
$wsErrores.PageSetup.PaperSize = 9
$wsErrores.PageSetup.Orientation = 1

# --- Sheet view / selection updates: "Datos" becomes the active/selected tab,
# --- "Errores" loses tabSelected, and the in-sheet selections move.
[void]$wsErrores.Range("C11").Select()
[void]$wsDatos.Activate()
[void]$wsDatos.Range("A38").Select()
